$wb = $excel.ActiveWorkbook

function Set-HandbackLink {
    param($ws, $row, $col, $url, $name)

    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $name
    $ws.Hyperlinks.Add($cell, $url, "", "", $name)
    # Match the existing "HyperLink" look (underline + FF6495ED) used by
    # the other linked cells in this workbook, rather than the
    # auto-generated theme hyperlink style that .Hyperlinks.Add() applies.
    $cell.Font.Name = "Calibri"
    $cell.Font.Underline = 2
    $cell.Font.Color = 0xED9564
}

function Update-HandbackSheet {
    param($ws, $mdUrl2, $mdUrl3, $xlfUrl2, $xlfUrl3, $xlfName2, $xlfName3, $handbackDateTime)

    # Column C (Status): "Ready for handoff" -> "Handed back: in sync with en-US"
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # Column F (Latest Target File) / Column G (Latest Handback File):
    # populated to mirror the source .md (col A) and handoff .xlf (col D)
    # files, since the handback is in sync with en-US (no new file produced).
    Set-HandbackLink $ws 2 6 $mdUrl2 "2f1e9fa4-5b61-4ab6-bbd4-2d02e34d2bed.md"
    Set-HandbackLink $ws 2 7 $xlfUrl2 $xlfName2
    Set-HandbackLink $ws 3 6 $mdUrl3 "f8cc793f-83df-4704-8c5e-8e3ac0b5ec84.md"
    Set-HandbackLink $ws 3 7 $xlfUrl3 $xlfName3

    # Column H (Latest Handback DateTime)
    $ws.Range("H2").Value = $handbackDateTime
    $ws.Range("H3").Value = $handbackDateTime
}

$mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/46795f18cebfc283cf0bcc9d9c68128b8b6054a3/e2e/2f1e9fa4-5b61-4ab6-bbd4-2d02e34d2bed.md"
$mdUrl3 = "https://github.com/OpenLocalizationTest/oltest/blob/46795f18cebfc283cf0bcc9d9c68128b8b6054a3/e2e/f8cc793f-83df-4704-8c5e-8e3ac0b5ec84.md"

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
Update-HandbackSheet $wsZh $mdUrl2 $mdUrl3 `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc72e4d718d292a95935aef0169f5e0b1d1c1d7b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/2f1e9fa4-5b61-4ab6-bbd4-2d02e34d2bed.f5b6a3af80fb18f2d58a3a837d372d82af50a741.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc72e4d718d292a95935aef0169f5e0b1d1c1d7b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/f8cc793f-83df-4704-8c5e-8e3ac0b5ec84.bd1213ba61a93e67b4bb3b5c5bb6ba008dd88ce8.zh-cn.xlf" `
    "2f1e9fa4-5b61-4ab6-bbd4-2d02e34d2bed.f5b6a3af80fb18f2d58a3a837d372d82af50a741.zh-cn.xlf" `
    "f8cc793f-83df-4704-8c5e-8e3ac0b5ec84.bd1213ba61a93e67b4bb3b5c5bb6ba008dd88ce8.zh-cn.xlf" `
    "2016-03-24 06:16:01"

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
Update-HandbackSheet $wsDe $mdUrl2 $mdUrl3 `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/65c2ff0135136443b67f66e0e75b3cc17a80c6d4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/2f1e9fa4-5b61-4ab6-bbd4-2d02e34d2bed.f5b6a3af80fb18f2d58a3a837d372d82af50a741.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/65c2ff0135136443b67f66e0e75b3cc17a80c6d4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/f8cc793f-83df-4704-8c5e-8e3ac0b5ec84.bd1213ba61a93e67b4bb3b5c5bb6ba008dd88ce8.de-de.xlf" `
    "2f1e9fa4-5b61-4ab6-bbd4-2d02e34d2bed.f5b6a3af80fb18f2d58a3a837d372d82af50a741.de-de.xlf" `
    "f8cc793f-83df-4704-8c5e-8e3ac0b5ec84.bd1213ba61a93e67b4bb3b5c5bb6ba008dd88ce8.de-de.xlf" `
    "2016-03-24 06:16:09"
